# Add a new entry (row 20) to the hours log on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data: date, hours, task, people(blank), paid
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A20").Value = 45776

$ws.Range("B20").Value = "2 hours"
$ws.Range("C20").Value = "update data and plots"
$ws.Range("E20").Value = "N"

# Match the selection left behind after the edit
$ws.Range("B18").Select()
